$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log entry in row 8: same "alleen" person as rows 4-7, new status text,
# same date as row 7 (23-10-18) with a begin time equal to row 7's end time
# (15:00) and a new end time of 16:20.
$ws.Range("A8").Value = "alleen"
$ws.Range("B8").Value = "parser werkend en geïnplementeerd en alle rest-protocollen werken nu(alleen nog niet alle foutmeldingen)"

$ws.Range("C8").Value = 43396
$ws.Range("C8").NumberFormat = "DD/MM/YY"

$ws.Range("D8").Value = 0.625
$ws.Range("D8").NumberFormat = "HH:MM:SS"

$ws.Range("E8").Value = 0.680555555555556
$ws.Range("E8").NumberFormat = "HH:MM:SS"

# Column B grew wider to fit the longer log text.
$ws.Columns("B").ColumnWidth = 89.67

# Selection moved to D10 after the edit.
[void]$ws.Range("D10").Select()
